# Team03Report.xlsx update — "update list of thigns to avoid and keep doing in
# spring report"
#
# 1. Sprint1 sheet: add new bullet detail rows under "Keep doing:" (row 16)
#    and "Avoid:" (originally row 20, renumbered to row 21 once the new rows
#    are inserted above it).
# 2. Burndown README sheet: scroll the view down so row 13 is at the top.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sprint1 sheet — retrospective notes
# ---------------------------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Activate()

# Keep the "Keep doing:" label's row (16) and give it supporting detail in
# column C, plus two brand-new detail-only rows below it.
$sprint1.Range("C16").Value = "Pre sprint planning meeting to discuss implementation pattern"
$sprint1.Range("C17").Value = "End of sprint code review"
$sprint1.Range("C18").Value = "List of common formatting in the code to keep style the same"

# Grab the formatting of the "Keep doing:" label cell so the relocated
# "Avoid:" label keeps its original look (bold/shaded style index 5) once we
# rebuild it a few rows further down.
$keepDoingLabel = $sprint1.Range("B16")

# The old "Avoid:" row (row 20) is removed outright; it reappears below with
# the new supporting rows interleaved.
$sprint1.Rows.Item(20).Delete()

$avoidLabel = $sprint1.Range("B21")
$keepDoingLabel.Copy()
$avoidLabel.PasteSpecial(-4122)  # xlPasteFormats
$avoidLabel.Value = "Avoid:"

$sprint1.Range("C21").Value = "Pre sprint meeting should've been earlier; people started working in different directions before we decided on a infrastructure. "
$sprint1.Range("C22").Value = "Mid sprint code review to make sure everyone is on the same page and catch early bugs"
$sprint1.Range("C23").Value = "Potential code refactor into different files to reduce file size. "

# Restore/refresh the sheet's selection to match the authored file (G16).
$sprint1.Range("G16").Select() | Out-Null

# ---------------------------------------------------------------------------
# Burndown README sheet — scroll position
# ---------------------------------------------------------------------------
$readme = $wb.Worksheets.Item("Burndown README")
$readme.Activate()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1

# Sprint1 remains the tab that is active/selected in the saved workbook, so
# switch back to it last.
$sprint1.Activate()
